# Update "想去人数" (F column) counts on both the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 819
    "F5"  = 46
    "F6"  = 12298
    "F7"  = 49
    "F10" = 436
    "F12" = 898
    "F13" = 13597
    "F14" = 13769
    "F19" = 1028
    "F20" = 103
    "F22" = 4851
    "F23" = 211
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
